# Refresh the "cryptos" price table (Wed Nov 22 04:07:32 UTC 2023 GitHub Actions run).
# Columns: A=rank index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Only B/C/D/E are ever touched; the workbook stores every one of these as
# literal text (t="inlineStr" in the source), so column D (which often looks
# numeric, e.g. "234.11") is forced to Text format before the write so Excel
# does not silently convert it to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.360.96'
$ws.Range("E2").Value = '  -2.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.978.77'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.11'
$ws.Range("E5").Value = '  -11.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.598'
$ws.Range("E6").Value = '  -3.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.28'
$ws.Range("E8").Value = '  -3.38%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.372'
$ws.Range("E9").Value = '  -3.51%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.89'
$ws.Range("E10").Value = '  +3.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0745'
$ws.Range("E11").Value = '  -3.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0984'
$ws.Range("E12").Value = '  -3.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.269.34'
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '13.99'
$ws.Range("E14").Value = '  -3.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.87'
$ws.Range("E15").Value = '  -4.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.750'
$ws.Range("E16").Value = '  -7.15%  '
$ws.Range("E17").Value = '  -4.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.975.06'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.364.65'
$ws.Range("E19").Value = '  -2.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.48'
$ws.Range("E20").Value = '  -3.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0802'
$ws.Range("E21").Value = '  -4.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.29'
$ws.Range("E22").Value = '  +1.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '221.33'
$ws.Range("E23").Value = '  -3.27%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +1.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.39'
$ws.Range("E26").Value = '  -11.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.58'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("E28").Value = '  -5.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.126'
$ws.Range("E29").Value = '  -2.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.74'
$ws.Range("E30").Value = '  -5.64%  '
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("E32").Value = '  -3.03%  '
$ws.Range("E33").Value = '  -6.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0603'
$ws.Range("E34").Value = '  -7.62%  '
$ws.Range("E35").Value = '  -6.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.31'
$ws.Range("E36").Value = '  -3.57%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.23'
$ws.Range("E39").Value = '  -2.85%  '
$ws.Range("E40").Value = '  +4.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.03'
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.456.15'
$ws.Range("E42").Value = '  +4.50%  '
$ws.Range("E43").Value = '  -3.87%  '
$ws.Range("E44").Value = '  -6.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.09'
$ws.Range("E45").Value = '  -10.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.01'
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '14.81'
$ws.Range("E47").Value = '  -5.79%  '
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.87'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.75'
$ws.Range("E50").Value = '  -4.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.161.05'
$ws.Range("E51").Value = '  -2.08%  '
